$wb = $excel.ActiveWorkbook

# ---- Sheet ALC: 66 cell updates ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2541.625
$ws.Range("J32").Value = 2476.1428
$ws.Range("L32").Value = 2476.1428
$ws.Range("N32").Value = -3128.1428
$ws.Range("H40").Value = 2329.9
$ws.Range("I40").Value = 1799
$ws.Range("J40").Value = 2388.889
$ws.Range("K40").Value = 1799
$ws.Range("L40").Value = 2388.889
$ws.Range("M40").Value = -1624
$ws.Range("N40").Value = -2738.889
$ws.Range("H53").Value = 3172.9
$ws.Range("I53").Value = 69.875
$ws.Range("J53").Value = 5241.5835
$ws.Range("K53").Value = 69.875
$ws.Range("L53").Value = 5241.5835
$ws.Range("M53").Value = 567.125
$ws.Range("N53").Value = -6515.5835
$ws.Range("H70").Value = 1460439.4
$ws.Range("J70").Value = 3199
$ws.Range("L70").Value = 9597
$ws.Range("N70").Value = -10137
$ws.Range("H73").Value = 1460439.4
$ws.Range("J73").Value = 3199
$ws.Range("L73").Value = 9597
$ws.Range("N73").Value = -11469
$ws.Range("H88").Value = 2080.8235
$ws.Range("I88").Value = 2149.7
$ws.Range("J88").Value = 1982.4286
$ws.Range("K88").Value = 2149.7
$ws.Range("L88").Value = 1982.4286
$ws.Range("M88").Value = -1743.7
$ws.Range("N88").Value = -2794.4286
$ws.Range("H91").Value = 2080.8235
$ws.Range("I91").Value = 2149.7
$ws.Range("J91").Value = 1982.4286
$ws.Range("K91").Value = 2149.7
$ws.Range("L91").Value = 1982.4286
$ws.Range("M91").Value = -745.6999999999998
$ws.Range("N91").Value = -4790.4286
$ws.Range("H111").Value = 13627.111
$ws.Range("I111").Value = 988
$ws.Range("K111").Value = 2964
$ws.Range("M111").Value = 103
$ws.Range("H112").Value = 1815.2122
$ws.Range("I112").Value = 1400
$ws.Range("J112").Value = 1828.1875
$ws.Range("K112").Value = 4200
$ws.Range("L112").Value = 5484.5625
$ws.Range("M112").Value = -3092
$ws.Range("N112").Value = -7700.5625
$ws.Range("H113").Value = 111114540
$ws.Range("J113").Value = 3985
$ws.Range("L113").Value = 3985
$ws.Range("N113").Value = -10493
$ws.Range("H116").Value = 53147440
$ws.Range("I116").Value = 50229830
$ws.Range("J116").Value = 55578780
$ws.Range("K116").Value = 50229830
$ws.Range("L116").Value = 55578780
$ws.Range("M116").Value = -50226388
$ws.Range("N116").Value = -55585664
$ws.Range("H141").Value = 2934.72
$ws.Range("I141").Value = 979.9524
$ws.Range("K141").Value = 2939.8572
$ws.Range("M141").Value = 2240.1428

# ---- Sheet ARM: 7 cell updates ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2149.18
$ws.Range("I32").Value = 2348.9778
$ws.Range("J32").Value = 351
$ws.Range("K32").Value = 2348.9778
$ws.Range("L32").Value = 351
$ws.Range("M32").Value = -2061.9778
$ws.Range("N32").Value = -925

# ---- Sheet BSM: 12 cell updates ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 11110.839
$ws.Range("I20").Value = 13797.762
$ws.Range("K20").Value = 13797.762
$ws.Range("M20").Value = -13550.762
$ws.Range("H94").Value = 1326.6666
$ws.Range("J94").Value = 5000
$ws.Range("L94").Value = 5000
$ws.Range("N94").Value = -5902
$ws.Range("H134").Value = 1668.3636
$ws.Range("I134").Value = 1594.1428
$ws.Range("K134").Value = 4782.428400000001
$ws.Range("M134").Value = -2247.428400000001

# ---- Sheet CRP: 25 cell updates ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 60.46154
$ws.Range("I7").Value = 68.77778000000001
$ws.Range("J7").Value = 41.75
$ws.Range("K7").Value = 68.77778000000001
$ws.Range("L7").Value = 41.75
$ws.Range("M7").Value = 44.22221999999999
$ws.Range("N7").Value = -267.75
$ws.Range("H80").Value = 49000
$ws.Range("I80").Value = 45000
$ws.Range("J80").Value = 50000
$ws.Range("K80").Value = 45000
$ws.Range("L80").Value = 50000
$ws.Range("M80").Value = -43877
$ws.Range("N80").Value = -52246
$ws.Range("H83").Value = 49000
$ws.Range("I83").Value = 45000
$ws.Range("J83").Value = 50000
$ws.Range("K83").Value = 135000
$ws.Range("L83").Value = 150000
$ws.Range("M83").Value = -129384
$ws.Range("N83").Value = -161232
$ws.Range("H134").Value = 2289.2173
$ws.Range("I134").Value = 1310.5333
$ws.Range("K134").Value = 3931.5999
$ws.Range("M134").Value = -1396.5999

# ---- Sheet CUL: 20 cell updates ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 366.66666
$ws.Range("I22").Value = 366.66666
$ws.Range("K22").Value = 1099.99998
$ws.Range("M22").Value = -930.9999800000001
$ws.Range("H26").Value = 217.85715
$ws.Range("J26").Value = 400
$ws.Range("L26").Value = 1200
$ws.Range("N26").Value = -1776
$ws.Range("H27").Value = 366.66666
$ws.Range("I27").Value = 366.66666
$ws.Range("K27").Value = 1099.99998
$ws.Range("M27").Value = -997.9999800000001
$ws.Range("H92").Value = 299.5
$ws.Range("I92").Value = 299.5
$ws.Range("K92").Value = 898.5
$ws.Range("M92").Value = 349.5
$ws.Range("H121").Value = 122130
$ws.Range("J121").Value = 309710
$ws.Range("L121").Value = 929130
$ws.Range("N121").Value = -931750

# ---- Sheet GSM: 15 cell updates ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 4281.5454
$ws.Range("I55").Value = 3378.875
$ws.Range("J55").Value = 6688.6665
$ws.Range("K55").Value = 3378.875
$ws.Range("L55").Value = 6688.6665
$ws.Range("M55").Value = -3051.875
$ws.Range("N55").Value = -7342.6665
$ws.Range("H123").Value = 89900
$ws.Range("J123").Value = 89900
$ws.Range("L123").Value = 89900
$ws.Range("N123").Value = -94800
$ws.Range("H131").Value = 41166.668
$ws.Range("J131").Value = 41166.668
$ws.Range("L131").Value = 41166.668
$ws.Range("N131").Value = -51246.668

# ---- Sheet LTW: 21 cell updates ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 8733.538
$ws.Range("I68").Value = 2749.5
$ws.Range("J68").Value = 18308
$ws.Range("K68").Value = 2749.5
$ws.Range("L68").Value = 18308
$ws.Range("M68").Value = -2000.5
$ws.Range("N68").Value = -19806
$ws.Range("H71").Value = 8733.538
$ws.Range("I71").Value = 2749.5
$ws.Range("J71").Value = 18308
$ws.Range("K71").Value = 13747.5
$ws.Range("L71").Value = 91540
$ws.Range("M71").Value = -10003.5
$ws.Range("N71").Value = -99028
$ws.Range("H132").Value = 5677.143
$ws.Range("I132").Value = 3609.5
$ws.Range("J132").Value = 6825.8335
$ws.Range("K132").Value = 10828.5
$ws.Range("L132").Value = 20477.5005
$ws.Range("M132").Value = -8298.5
$ws.Range("N132").Value = -25537.5005

# ---- Sheet WVR: 12 cell updates ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 30000
$ws.Range("I63").Value = 30000
$ws.Range("K63").Value = 30000
$ws.Range("M63").Value = -29376
$ws.Range("H66").Value = 30000
$ws.Range("I66").Value = 30000
$ws.Range("K66").Value = 90000
$ws.Range("M66").Value = -86880
$ws.Range("H132").Value = 4223.081
$ws.Range("I132").Value = 4041.7856
$ws.Range("K132").Value = 12125.3568
$ws.Range("M132").Value = -9595.356800000001

